$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B21").Value = "5-75 Manufacturing `n5-30 Services and others"
$ws.Range("D21").Value = "300,000-15Million RM Manufacturing `n300,000-3Million RM Services & others"
$ws.Range("B22").Value = "75-200 Manufacturing `n30-75 Services  and others"
$ws.Range("D22").Value = "15-50Million RM Manufacturing `n3-20Million RM Services & others"
$ws.Range("B23").Value = ">200 Manufacturing `n>75 Services  and others"
$ws.Range("D23").Value = ">50Million RM Manufacturing `n>20Million RM Services & others"
